# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") -- theme1.xml becomes the Integral
# theme and theme2.xml becomes the (generic) Office Theme.
#
# ppt/theme/theme2.xml is the theme actually wired to the slide master /
# slides (ppt/_rels/presentation.xml.rels + slideMaster1.xml.rels both
# point at theme2.xml), so this is the part of the swap that is visible
# across every slide. The only substantive content difference between the
# two theme parts is their 12-colour <a:clrScheme> (fonts/format scheme are
# byte-identical), so recolour the live theme's scheme to what is
# currently theme1.xml's ("Office Theme") palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Target palette = the current "Office Theme" colours (theme1.xml),
# ordered dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (ThemeColorScheme
# index order 1-12).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
